# Applies the row-data correction for rows 55-61 of sheet "Artfynd":
# species records for rows 55-59 are re-ordered (cyclic shift) and the
# Ost/Nord (Q/R) coordinates for rows 55-61 are rounded to whole metres;
# the Starttid/Sluttid ("00:00") cells (Z/AB) are dropped for rows 55-61.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 55 ---
$ws.Range("A55").Value = 111898191
$ws.Range("B55").Value = 90332
$ws.Range("D55").Value = "LC"
$ws.Range("E55").Value = 4769
$ws.Range("F55").Value = "Svavelriska"
$ws.Range("G55").Value = "Lactarius scrobiculatus"
$ws.Range("H55").Value = "(Scop.:Fr.) Fr."
$ws.Range("I55").Value = "2"
$ws.Range("J55").Value = "fruktkroppar"
$ws.Range("L55").ClearContents()
$ws.Range("Q55").Value = 650135
$ws.Range("R55").Value = 6654003
$ws.Range("Z55").ClearContents()
$ws.Range("AB55").ClearContents()
$ws.Range("AC55").ClearContents()
$ws.Range("AH55").Value = "Ängsbarrskog"
$ws.Range("AI55").Value = "Ungskog"

# --- Row 56 ---
$ws.Range("A56").Value = 111898889
$ws.Range("B56").Value = 98535
$ws.Range("E56").Value = 222498
$ws.Range("F56").Value = "Blåsippa"
$ws.Range("G56").Value = "Hepatica nobilis"
$ws.Range("H56").Value = "Schreb."
$ws.Range("I56").ClearContents()
$ws.Range("J56").ClearContents()
$ws.Range("K56").Value = "fullt utvecklade blad"
$ws.Range("L56").ClearContents()
$ws.Range("Q56").Value = 650135
$ws.Range("R56").Value = 6654003
$ws.Range("Z56").ClearContents()
$ws.Range("AB56").ClearContents()

# --- Row 57 ---
$ws.Range("A57").Value = 111898507
$ws.Range("B57").Value = 89845
$ws.Range("D57").Value = "VU"
$ws.Range("E57").Value = 1209
$ws.Range("F57").Value = "Rynkskinn"
$ws.Range("G57").Value = "Phlebia centrifuga"
$ws.Range("H57").Value = "P.Karst."
$ws.Range("Q57").Value = 650087
$ws.Range("R57").Value = 6654015
$ws.Range("Z57").ClearContents()
$ws.Range("AB57").ClearContents()

# --- Row 58 ---
$ws.Range("A58").Value = 111898660
$ws.Range("B58").Value = 100532
$ws.Range("D58").Value = "CR"
$ws.Range("E58").Value = 223246
$ws.Range("F58").Value = "Skogsalm"
$ws.Range("G58").Value = "Ulmus glabra"
$ws.Range("H58").Value = "Huds."
$ws.Range("L58").ClearContents()
$ws.Range("Q58").Value = 650054
$ws.Range("R58").Value = 6654018
$ws.Range("Z58").ClearContents()
$ws.Range("AB58").ClearContents()
$ws.Range("AC58").Value = "Stammens omkrets i brösthöjd: 64 cm"
$ws.Range("AJ58").ClearContents()
$ws.Range("AK58").ClearContents()
$ws.Range("AM58").ClearContents()
$ws.Range("AO58").ClearContents()

# --- Row 59 ---
$ws.Range("A59").Value = 111898336
$ws.Range("B59").Value = 89405
$ws.Range("D59").Value = "NT"
$ws.Range("E59").Value = 1202
$ws.Range("F59").Value = "Ullticka"
$ws.Range("G59").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H59").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("K59").ClearContents()
$ws.Range("L59").ClearContents()
$ws.Range("Q59").Value = 650105
$ws.Range("R59").Value = 6654011
$ws.Range("Z59").ClearContents()
$ws.Range("AB59").ClearContents()
$ws.Range("AH59").Value = "Ängsblandskog"
$ws.Range("AI59").ClearContents()
$ws.Range("AJ59").Value = "gran"
$ws.Range("AK59").Value = "Picea abies"
$ws.Range("AM59").Value = "Liggande död trädstam, utan markontakt"
$ws.Range("AO59").Value = "Horizontal, dead without ground contact # Picea abies"

# --- Row 60 ---
$ws.Range("Q60").Value = 650033
$ws.Range("R60").Value = 6654279
$ws.Range("Z60").ClearContents()
$ws.Range("AB60").ClearContents()

# --- Row 61 ---
$ws.Range("Q61").Value = 650027
$ws.Range("R61").Value = 6654299
$ws.Range("Z61").ClearContents()
$ws.Range("AB61").ClearContents()
